$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "Neuse River Brewing" entry as row 28 (new shared strings are
# created automatically for the new text values).
$ws.Range("A28").Value = "Neuse River Brewing"
$ws.Range("B28").Value = "Raleigh"
$ws.Range("C28").Value = "Neuse River Burger"
$ws.Range("D28").Value = "Brassiere/Burgers"
$ws.Range("E28").Value = 35.804566959442603
$ws.Range("F28").Value = -78.632520307935593

# Match the author's final UI state: scrolled/selected on C11.
$ws.Range("C11").Select()
